$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Starting layout (before edit):
#   row 6  : header
#   row 7  : item #1 - CETAL 1000MG 15 TABS
#   row 8  : item #2 - COLD FREE 20 TAB.
#   row 9  : totals
#   row 10 : footer (timestamp / page / credit)
#
# Target layout (after edit):
#   row 6  : header
#   row 7  : item #1 - ANTODINE 20MG 30 F.C.TAB   (NEW)
#   row 8  : item #2 - CETAL 1000MG 15 TABS       (was row 7)
#   row 9  : item #3 - COLD FREE 20 TAB.          (was row 8)
#   row 10 : item #4 - LIBRAX 30 SUGAR COATED TAB (NEW)
#   row 11 : totals                               (was row 9, new sum)
#   row 12 : footer                                (was row 10, new timestamp)
# ---------------------------------------------------------------------------

# Insert a fresh row above the first item row - pushes CETAL/COLD FREE/totals/footer
# down by one (rows 7..10 -> 8..11).
$ws.Rows("7:7").Insert()

# Insert a fresh row above what is now the totals row (row 10) - pushes
# totals/footer down by one more (rows 10..11 -> 11..12), leaving a blank
# row 10 for the new LIBRAX line.
$ws.Rows("10:10").Insert()

# --- copy the per-cell formatting of the surviving item rows onto the two new rows ---
# Row 7 (new) should look like row 8 (CETAL - an "odd" item row, ht 25.5)
$ws.Range("A8:Q8").Copy()
$ws.Range("A7:Q7").PasteSpecial(-4122)

# Row 10 (new) should look like row 9 (COLD FREE - an "even" item row, ht 24.75)
$ws.Range("A9:Q9").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- row heights (alternating 25.5 / 24.75 by row position) ---
$ws.Rows("7:7").RowHeight = 25.5
$ws.Rows("8:8").RowHeight = 24.75
$ws.Rows("9:9").RowHeight = 25.5
$ws.Rows("10:10").RowHeight = 24.75

# --- merged cells for the two new rows (same pattern as the other item rows) ---
$ws.Range("A7:B7").Merge()
$ws.Range("C7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()
$ws.Range("N7:O7").Merge()

$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

# --- values for the new ANTODINE row (row 7) ---
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "ANTODINE 20MG 30 F.C.TAB"
$ws.Range("H7").Value = "1:0"
$ws.Range("L7").Value = "1"
$ws.Range("N7").Value = "60.00"
$ws.Range("P7").Value = "19.8000"
$ws.Range("Q7").Value = "0:1"

# --- renumber the existing item rows (now shifted to 8 and 9) ---
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 3

# --- values for the new LIBRAX row (row 10) ---
$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "LIBRAX 30 SUGAR COATED TAB"
$ws.Range("H10").Value = "3:2"
$ws.Range("L10").Value = "1"
$ws.Range("N10").Value = "48.00"
$ws.Range("P10").Value = "15.8400"
$ws.Range("Q10").Value = "0:1"

# --- updated total (row 11, was row 9) ---
$ws.Range("P11").Value = 74.5

# --- updated footer timestamp (row 12, was row 10) ---
$ws.Range("A12").Value = "Saturday, 9 August, 2025 9:31 AM"
